$wb = $excel.ActiveWorkbook
$win = $excel.ActiveWindow
$win.ScrollWorkbookTabs(1)
Write-Output "done"
